{"js": "// The placeholder text \"${schule}\" becomes \"${schule_nametype}\"\n// (Word split the edit into three runs while the author typed\n// \"_nametype\" just before the closing brace; the runs all share the\n// same formatting, so the visible/semantic result is a single merged\n// string - which is what we reproduce here.)\nconst body = context.document.body;\n\nconst results = body.search(\"${schule}\", { matchCase: true, matchWildcards: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"${schule_nametype}\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# The placeholder text \"${schule}\" becomes \"${schule_nametype}\"\n# (Word split the edit into three runs while the author typed\n# \"_nametype\" just before the closing brace; the runs all share the\n# same formatting, so the visible/semantic result is a single merged\n# string - which is what we reproduce here.)\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"`${schule}\"\n$find.MatchWildcards = $false\n$find.MatchCase = $true\n$find.Execute()\n\nif ($find.Found) {\n    $find.Parent.Text = \"`${schule_nametype}\"\n}\n"}
